$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 663, shifting existing rows 663-726 down to 666-729
$ws.Rows("663:665").Insert()

# New row 663: Papa / Asterix / 1a (guarda) / Región de Los Lagos
$ws.Range("A663").Value = 7
$ws.Range("B663").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C663").Value = "Ñuble"
$ws.Range("D663").Value = 45194
$ws.Range("E663").Value = 16
$ws.Range("F663").Value = 100114001
$ws.Range("G663").Value = "Papa"
$ws.Range("H663").Value = "Asterix"
$ws.Range("I663").Value = "1a (guarda)"
$ws.Range("J663").Value = 150
$ws.Range("K663").Value = 27000
$ws.Range("L663").Value = 27000
$ws.Range("M663").Value = 27000
$ws.Range("N663").Value = "$/saco 25 kilos"
$ws.Range("O663").Value = "Región de Los Lagos"
$ws.Range("P663").Value = 1080
$ws.Range("Q663").Value = 25
$ws.Range("R663").Value = "Hortaliza"

# New row 664: Papa / Asterix / 2a (guarda) / Región de Los Lagos
$ws.Range("A664").Value = 7
$ws.Range("B664").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C664").Value = "Ñuble"
$ws.Range("D664").Value = 45194
$ws.Range("E664").Value = 16
$ws.Range("F664").Value = 100114001
$ws.Range("G664").Value = "Papa"
$ws.Range("H664").Value = "Asterix"
$ws.Range("I664").Value = "2a (guarda)"
$ws.Range("J664").Value = 150
$ws.Range("K664").Value = 25000
$ws.Range("L664").Value = 25000
$ws.Range("M664").Value = 25000
$ws.Range("N664").Value = "$/saco 25 kilos"
$ws.Range("O664").Value = "Región de Los Lagos"
$ws.Range("P664").Value = 1000
$ws.Range("Q664").Value = 25
$ws.Range("R664").Value = "Hortaliza"

# New row 665: Papa / Rodeo / 1a (guarda) / Región de Los Lagos
$ws.Range("A665").Value = 7
$ws.Range("B665").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C665").Value = "Ñuble"
$ws.Range("D665").Value = 45194
$ws.Range("E665").Value = 16
$ws.Range("F665").Value = 100114001
$ws.Range("G665").Value = "Papa"
$ws.Range("H665").Value = "Rodeo"
$ws.Range("I665").Value = "1a (guarda)"
$ws.Range("J665").Value = 150
$ws.Range("K665").Value = 27000
$ws.Range("L665").Value = 27000
$ws.Range("M665").Value = 27000
$ws.Range("N665").Value = "$/saco 25 kilos"
$ws.Range("O665").Value = "Región de Los Lagos"
$ws.Range("P665").Value = 1080
$ws.Range("Q665").Value = 25
$ws.Range("R665").Value = "Hortaliza"
